$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newB2 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Not Applicable']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newB2

# Row 2 height shrinks now that the Cohort line has been removed (matches Excel's auto layout)
$ws.Rows.Item(2).RowHeight = 290

# Update the view: top-left cell and current selection move to B2
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select()
